# Updates cryptos list data (prices in column D, 1h volume % in column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.529.06'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '2.490.15'
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '''313.42'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").Value = '''93.24'
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("D7").Value = '''0.544'
$ws.Range("E7").Value = '  -1.52%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = '''0.500'
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("D10").Value = '''32.63'
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("D11").Value = '''0.0786'
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("E12").Value = '  +1.99%  '
$ws.Range("D13").Value = '2.876.19'
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("D14").Value = '''6.86'
$ws.Range("E14").Value = '  -2.61%  '
$ws.Range("D15").Value = '''16.23'
$ws.Range("E15").Value = '  +10.50%  '
$ws.Range("D16").Value = '2.505.00'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("D17").Value = '''0.757'
$ws.Range("E17").Value = '  -3.74%  '
$ws.Range("D18").Value = '41.566.11'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").Value = '''6.35'
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").Value = '0.0₃0930'
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").Value = '''71.32'
$ws.Range("E21").Value = '  +4.70%  '
$ws.Range("D22").Value = '''11.25'
$ws.Range("E22").Value = '  -2.47%  '
$ws.Range("D23").Value = '''236.13'
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("D24").Value = '''2.71'
$ws.Range("E24").Value = '  -3.03%  '
$ws.Range("E25").Value = '  -0.50%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").Value = '''25.37'
$ws.Range("E27").Value = '  +3.71%  '
$ws.Range("E28").Value = '  -0.60%  '
$ws.Range("D29").Value = '''9.69'
$ws.Range("E29").Value = '  -0.21%  '
$ws.Range("D30").Value = '''36.23'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").Value = '''157.81'
$ws.Range("E31").Value = '  +2.73%  '
$ws.Range("D32").Value = '''5.45'
$ws.Range("E32").Value = '  -2.68%  '
$ws.Range("D33").Value = '''2.57'
$ws.Range("E33").Value = '  -1.22%  '
$ws.Range("D34").Value = '''0.0758'
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D35").Value = '''17.86'
$ws.Range("E35").Value = '  +5.37%  '
$ws.Range("E36").Value = '  -6.27%  '
$ws.Range("D37").Value = '''2.95'
$ws.Range("E37").Value = '  -2.46%  '
$ws.Range("D38").Value = '''1.85'
$ws.Range("E38").Value = '  -2.25%  '
$ws.Range("E39").Value = '  +1.09%  '
$ws.Range("E40").Value = '  -0.42%  '
$ws.Range("D41").Value = '''4.13'
$ws.Range("E41").Value = '  -4.01%  '
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").Value = '''19.87'
$ws.Range("E43").Value = '  -6.85%  '
$ws.Range("D44").Value = '1.966.25'
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("E45").Value = '  -0.43%  '
$ws.Range("E46").Value = '  -2.78%  '
$ws.Range("E47").Value = '  +1.59%  '
$ws.Range("D48").Value = '2.729.00'
$ws.Range("E48").Value = '  +0.69%  '
$ws.Range("D49").Value = '''96.70'
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").Value = '''67.94'
$ws.Range("E50").Value = '  -3.11%  '
$ws.Range("D51").Value = '''73.83'
$ws.Range("E51").Value = '  -3.00%  '
